# Generate Report for Handoff
# Updates the "c0d0d805-ce1f-4c99-978b-71f547013d2b" row to reflect that the
# handoff/report generation has run again: Status becomes "Ready for handoff",
# the handoff timestamps are refreshed, and an Error Detail message is recorded
# because the handback file available is not the latest version.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$zhHandoffDatetime = "2016-10-20 00:23:17"
$deHandoffDatetime = "2016-10-20 00:23:28"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba9e33a27629a712daa8a033b37f7904c06f0154/e2e/c0d0d805-ce1f-4c99-978b-71f547013d2b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/655df9721bf41a9353aae3a55f9d49feb46a3a43/e2e/c0d0d805-ce1f-4c99-978b-71f547013d2b.md."

# ---- Overview sheet: row for c0d0d805 (row 3) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $deHandoffDatetime

# ---- zh-cn sheet: row for c0d0d805 (row 3) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("H3").Value = $zhHandoffDatetime
$zhcn.Range("P3").Value = $errorDetail
# widen the Error Detail column to fit the longer message
$zhcn.Range("P1").ColumnWidth = 39.183673469387755

# ---- de-de sheet: row for c0d0d805 (row 3) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("H3").Value = $deHandoffDatetime
$dede.Range("P3").Value = $errorDetail
# widen the Error Detail column to fit the longer message
$dede.Range("P1").ColumnWidth = 39.183673469387755
